$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","E","F","G","I","J","K","N")

$data = @{
  2 = @(0.5844107609055129,0.1603393385757954,0.3830702605408334,2.974491789485754,0.002507347822420498,1.043310562392463,0.1098224562254195,0.7315323515092871,2.074898677527742)
  3 = @(0.5445341824638774,0.1504787813262283,0.3653107241800342,2.934351508913011,0.0025110555045518,1.046327220668587,0.1098414624361794,0.682934081753416,2.095182165799457)
  4 = @(0.5203479662677921,0.1445238304920196,0.354612804271305,2.911429858977712,0.002513450718006661,1.048709053007826,0.1099139044731885,0.6534889827458414,2.108332443790353)
  5 = @(0.5105668322974566,0.1421219845714745,0.3503050837856279,2.902521718270933,0.00251445672724318,1.049812601953221,0.1099586680568585,0.6415888239766332,2.113866116251359)
  6 = @(0.5089472082655675,0.1417246555904654,0.3495929122204018,2.901068626732965,0.002514625585666616,1.050003867346689,0.1099670205214665,0.6396187821155763,2.1147955310844)
  7 = @(0.5202157510607321,0.1444913379837658,0.3545544994436653,2.911307970433171,0.002513464164135429,1.048723397934594,0.1099144465032644,0.6533280926603595,2.108406365267534)
  8 = @(0.5705993992849017,0.1569186683289558,0.3769038088499599,2.960292810033962,0.002508601657330624,1.04424064424353,0.1098163725891652,0.7146936487876019,2.081747809907078)
  9 = @(0.6717747354526296,0.1820869130135065,0.4223792506681718,3.070093231799206,0.002500003543264998,1.039663658355934,0.1101081405951803,0.8381784187685639,2.035003226755862)
  10 = @(0.7475738205519065,0.2010792749191523,0.4568128779563239,3.159233035162544,0.002494251648185358,1.038887095219764,0.1106204398331698,0.9308558378941711,2.004044668464076)
  11 = @(0.7823796787389483,0.2098315329916716,0.47270386642451,3.201645087382246,0.002491756347313045,1.039099258019569,0.1109188213614303,0.9734494934932059,1.990698353234315)
  12 = @(0.7956066114278144,0.2131621943268556,0.4787542641138742,3.217974636691906,0.002490828778619148,1.039261208412881,0.1110412542488035,0.9896414955020703,1.98575067536953)
  13 = @(0.7927558765242111,0.2124441456774377,0.4774497397508668,3.214445790408803,0.002491027777097585,1.039222695077775,0.1110144654286813,0.9861514648217167,1.986811514265725)
  14 = @(0.783466930011997,0.2101052193138457,0.4732009773196921,3.202983130127109,0.002491679688413104,1.039110944472021,0.1109287044781908,0.9747803609043331,1.990289172904653)
  15 = @(0.7777832666262157,0.2086746970218201,0.4706027681454401,3.195996993279465,0.002492081259739396,1.039053130899987,0.110877404497451,0.9678234060437205,1.992433189842416)
  16 = @(0.7453057093071891,0.2005095748631049,0.4557789475206704,3.156498883606247,0.00249441715459664,1.038884633374138,0.1106022586193234,0.92808100283969,2.00493173671827)
  17 = @(0.7254649196460718,0.195529493714389,0.4467432762646268,3.132745860659583,0.002495881144680798,1.038926297981718,0.1104502326447161,0.9038117306001539,2.012788149965616)
  18 = @(0.7140835598649744,0.1926756638951588,0.4415675480207284,3.119258923666081,0.002496734612643546,1.039003464596419,0.110368938012428,0.8898935427130823,2.017376291039213)
  19 = @(0.7102352704293367,0.1917112177734168,0.4398187967173612,3.114722525933161,0.002497025546340538,1.039038720301662,0.1103424673125133,0.8851880924278532,2.01894165638361)
  20 = @(0.7275738446947742,0.1960585353883175,0.4477029278079101,3.135256271431018,0.002495724119275855,1.038916354301527,0.1104657795698074,0.906391005751459,2.011944642633708)
  21 = @(0.7861940506535916,0.2107917728608868,0.4744480494427563,3.206342679691943,0.002491487736475489,1.039141551069683,0.1109536379077483,0.9781186256931562,1.9892648132316)
  22 = @(0.8247780307439143,0.2205163080210184,0.4921189539396096,3.254370454385622,0.002488820089259308,1.039764552017004,0.1113275407195573,1.025362417502436,1.975061878597572)
  23 = @(0.8041601019435518,0.2153173382806415,0.4826700810072992,3.22859316145491,0.002490234644163029,1.039388406429865,0.1111229280501007,1.000113969840953,1.982585447906459)
  24 = @(0.7266203206446278,0.1958193267405761,0.4472690101831347,3.134120788394114,0.002495795073597106,1.038920684093867,0.1104587317863235,0.9052248078882883,2.0123257699383)
  25 = @(0.6441481385176644,0.1751912735393262,0.4098986244266314,3.038908726877082,0.002502229868535008,1.040448967481709,0.1099770945790226,0.8044320492000736,2.04705535933752)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + $r).Value = $vals[$i]
  }
}
